# Updated symbol list on Sat Dec 17 05:49:23 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D are stored as plain text (not numbers) in this
# sheet. Assigning a numeric-looking string via .Value auto-converts the
# cell to a Number, so force a text number format first and restore a
# "clean" (unstyled) look afterwards by copying the style of an untouched
# text cell in the same column.
$refStyle = $ws.Range("D17").Style
function Set-TextValue($rangeAddress, $text) {
    $r = $ws.Range($rangeAddress)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = $refStyle
}

# --- Simple price (column D) refreshes ---
Set-TextValue "D2"  "229.49"
Set-TextValue "D3"  "22.61"
Set-TextValue "D4"  "5.274"
Set-TextValue "D5"  "0.05576"
Set-TextValue "D6"  "3.381"
Set-TextValue "D7"  "6.471"

# --- Row 8 / Row 9 swap (MXToken <-> FTXToken changed ranking order) ---
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D8" "1.045"
$ws.Range("E8").Value = "7FTXTokenFTT"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.7813"
$ws.Range("E9").Value = "8MXTokenMX"

# --- More simple price refreshes ---
Set-TextValue "D10" "0.1379"
Set-TextValue "D11" "0.07343"
Set-TextValue "D12" "0.03169"
Set-TextValue "D13" "0.02969"
Set-TextValue "D14" "0.09271"
Set-TextValue "D15" "0.001658"
Set-TextValue "D16" "3.264"

# --- Rows 18-24: ranking shuffled by one position, each row takes on the
#     coin that used to be one row below it; row 24 gets "One" back with a
#     refreshed price and a "Bestin24h" marker ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.006223"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D19" "0.005232"
$ws.Range("E19").Value = "18HotbitTokenHTB"

$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D20" "0.001061"
$ws.Range("E20").Value = "19BitKanKAN"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.0001501"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.945"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D23" "2.146"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D24" "0.01164"
$ws.Range("E24").Value = "23OneONEBestin24h"

# --- Row 27: price refresh + drop the stray "Bestin24h" tag from E27 ---
Set-TextValue "D27" "0.0005004"
$ws.Range("E27").Value = "26UpBotsUBXT"

# --- More simple price refreshes ---
Set-TextValue "D40" "0.03996"

# --- Rows 42 / 43 swap (CEJI <-> BKEXToken changed ranking order) ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003406"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1037"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# --- Remaining simple price refreshes ---
Set-TextValue "D44" "0.01000"
Set-TextValue "D45" "0.00005445"
Set-TextValue "D47" "0.7859"
Set-TextValue "D48" "0.04238"
